$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2020-11-21")

# New attendance rows for "dishant" (training set addition) following the
# same pattern/columns as the existing rows in this sheet (A:I).
$rows = @(
    @{ Row = 7; Time = "02:02:54"; SpO2 = 99.88070376432579; HR = 57.52758683819931 },
    @{ Row = 8; Time = "02:04:38"; SpO2 = 97.03971006047878; HR = 89.00981767428932 },
    @{ Row = 9; Time = "02:06:08"; SpO2 = 98.84141964022119; HR = 82.47545924642802 }
)

foreach ($r in $rows) {
    $i = $r.Row

    $ws.Range("A$i").Value = 1
    $ws.Range("B$i").Value = "sachin"
    $ws.Range("C$i").Value = "301/Sanskruti-1,Andheri, Mumbai"
    $ws.Range("D$i").Value = "Software Engineer"
    $ws.Range("E$i").Value = $r.Time
    $ws.Range("F$i").Value = $r.SpO2
    $ws.Range("G$i").Value = $r.HR
    $ws.Range("H$i").Value = "NA"
    $ws.Range("I$i").Value = "NA"

    # Match the bold / centered / thin-bordered look used by column A on the
    # existing rows (A2:A6) for the "Sr. No" column.
    $cell = $ws.Range("A$i")
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
}

Write-Output "added rows 7-9 to 2020-11-21"
